$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: count (31 -> 34) for columns B..T
$ws.Range("B2:T2").Value = 34

# Row 3: mean
$ws.Range("B3").Value = 3.058823529411764
$ws.Range("C3").Value = 3.352941176470588
$ws.Range("D3").Value = 2.794117647058823
$ws.Range("E3").Value = 2.764705882352941
$ws.Range("F3").Value = 2.794117647058823
$ws.Range("G3").Value = 3.294117647058823
$ws.Range("H3").Value = 2.352941176470588
$ws.Range("I3").Value = 2.617647058823529
$ws.Range("J3").Value = 2.705882352941177
$ws.Range("K3").Value = 3.264705882352941
$ws.Range("L3").Value = 3.352941176470588
$ws.Range("M3").Value = 3.794117647058823
$ws.Range("N3").Value = 3.205882352941177
$ws.Range("O3").Value = 3.529411764705882
$ws.Range("P3").Value = 4.088235294117647
$ws.Range("Q3").Value = 3.411764705882353
$ws.Range("R3").Value = 2.823529411764706
$ws.Range("S3").Value = 3.264705882352941
$ws.Range("T3").Value = 2.676470588235294

# Row 4: std
$ws.Range("B4").Value = 1.347077115755897
$ws.Range("C4").Value = 1.29993829553272
$ws.Range("D4").Value = 1.628894619890878
$ws.Range("E4").Value = 1.538571228778491
$ws.Range("F4").Value = 1.430816926321757
$ws.Range("G4").Value = 1.487919266397041
$ws.Range("H4").Value = 1.411690431795447
$ws.Range("I4").Value = 1.51786744105915
$ws.Range("J4").Value = 1.467411899471122
$ws.Range("K4").Value = 1.377499348939224
$ws.Range("L4").Value = 1.495090003192804
$ws.Range("M4").Value = 1.174976769975549
$ws.Range("N4").Value = 1.365803388057981
$ws.Range("O4").Value = 1.236693884801685
$ws.Range("P4").Value = 1.13798458366358
$ws.Range("Q4").Value = 1.233807780786253
$ws.Range("R4").Value = 1.566129919052426
$ws.Range("S4").Value = 1.377499348939224
$ws.Range("T4").Value = 1.570959718930652

# Row 6: 25% quartile (only some cells change)
$ws.Range("E6").Value = 1.25
$ws.Range("J6").Value = 1.25
$ws.Range("L6").Value = 2
$ws.Range("O6").Value = 3
$ws.Range("R6").Value = 1.25

# Row 7: 50% quartile (median, only some cells change)
$ws.Range("C7").Value = 3.5
$ws.Range("L7").Value = 3.5
$ws.Range("R7").Value = 2
$ws.Range("T7").Value = 2.5

# Row 8: 75% quartile (only some cells change)
$ws.Range("G8").Value = 4.75
$ws.Range("I8").Value = 3.75
$ws.Range("K8").Value = 4
$ws.Range("O8").Value = 4.75
$ws.Range("R8").Value = 4
